# Chapter 02 "Entering Data" exercise file - fill in the DataEntry sheet
# with the worked example content (expenses table, autofill quarters,
# date/time entry demo) and the author's review comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# ---- Row 2: Sales ------------------------------------------------------
$ws.Range("A2").Value = "Sales"
$ws.Range("B2").Value = 120
$ws.Range("C2").Value = 160
$ws.Range("D2").Value = 190
$ws.Range("E2").Value = 220
$ws.Range("F2").Value = 240
$ws.Range("G2").Value = 290

# ---- Row 3: Expenses ----------------------------------------------------
$ws.Range("A3").Value = "Expenses"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 160
$ws.Range("E3").Value = 170
$ws.Range("F3").Value = 200
$ws.Range("G3").Value = 210

# ---- Row 4: Profits label ------------------------------------------------
$ws.Range("A4").Value = "Profits"

# ---- Row 1: month headers (autofilled January..June) --------------------
$ws.Range("B1").Value = "January"
$ws.Range("C1").Value = "February"
$ws.Range("D1").Value = "March"
$ws.Range("E1").Value = "April"
$ws.Range("F1").Value = "May"
$ws.Range("G1").Value = "June"

# ---- Row 6: Quarter autofill demo ---------------------------------------
$ws.Range("B6").Value = "Quarter1"
$ws.Range("C6").Value = "Quarter2"
$ws.Range("D6").Value = "Quarter3"
$ws.Range("E6").Value = "Quarter4"
$ws.Range("F6").Value = "Quarter1"
$ws.Range("G6").Value = "Quarter2"

# ---- Row 9: Date / Time headers ------------------------------------------
$ws.Range("B9").Value = "Date"
$ws.Range("D9").Value = "Time"

# ---- Row 10: date entered with "/" shown with long date format, time with h:mm
$ws.Range("B10").Value = 44876
$ws.Range("B10").NumberFormat = "mm-dd-yy"
$ws.Range("D10").Value = 0.49236111111111108
$ws.Range("D10").NumberFormat = "h:mm"

# ---- Row 11: same date/time, AM/PM time format ---------------------------
$ws.Range("B11").Value = 44876
$ws.Range("B11").NumberFormat = "mm-dd-yy"
$ws.Range("D11").Value = 0.49236111111111108
$ws.Range("D11").NumberFormat = "h:mm AM/PM"

# ---- Row 12: invalid date entered as text ---------------------------------
$ws.Range("B12").Value = "32-11-22"

# ---- Row 13: formula - 100 days after B10 ----------------------------------
$ws.Range("B13").Formula = "=B10+100"
$ws.Range("B13").NumberFormat = "mm-dd-yy"

# ---- Column B a bit wider, to show invalid dates clearly -------------------
$ws.Columns.Item(2).ColumnWidth = 9.86

# ---- View: zoomed in to 200%, scrolled down, D12 selected -------------------
$excel.ActiveWindow.Zoom = 200
$ws.Range("D12").Select()

# ---- Review comments left by Islam, Md Muntaha EX1 --------------------------
$author = "Islam, Md Muntaha EX1:"

$c = $ws.Range("G1").AddComment($author + "`nAutofilled.")
$c = $ws.Range("A2").AddComment($author + "`nStrings alined to left.")
$c = $ws.Range("B2").AddComment($author + "`nNumbers alined to right.")
$c = $ws.Range("E6").AddComment($author + "`nThis autofill stops at 4 and then restarts.")
$c = $ws.Range("B9").AddComment($author + "`nKeep date columns a bit wider to find out if a date is invalid.")
$c = $ws.Range("B11").AddComment($author + "`nEven after I put slashes, it changes to dash auto as that’s my pc's date system.")
$c = $ws.Range("D11").AddComment($author + "`nam -> AM")
$c = $ws.Range("B12").AddComment($author + "`nAs it's invalid, it's alined to left, as a string.")
$c = $ws.Range("B13").AddComment($author + "`n100 days after day on B10.")
